$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 24428.572
$ws.Range("I21").Value = 19000
$ws.Range("J21").Value = 28500
$ws.Range("K21").Value = 19000
$ws.Range("L21").Value = 28500
$ws.Range("M21").Value = -18532
$ws.Range("N21").Value = -29436
$ws.Range("H23").Value = 24428.572
$ws.Range("I23").Value = 19000
$ws.Range("J23").Value = 28500
$ws.Range("K23").Value = 19000
$ws.Range("L23").Value = 28500
$ws.Range("M23").Value = -18766
$ws.Range("N23").Value = -28968
$ws.Range("H28").Value = 577197.5600000001
$ws.Range("I28").Value = 757165.5600000001
$ws.Range("K28").Value = 757165.5600000001
$ws.Range("M28").Value = -756680.5600000001
$ws.Range("H46").Value = 3845.9092
$ws.Range("I46").Value = 4502.4287
$ws.Range("J46").Value = 2697
$ws.Range("K46").Value = 13507.2861
$ws.Range("L46").Value = 8091
$ws.Range("M46").Value = -13388.2861
$ws.Range("N46").Value = -8329
$ws.Range("H55").Value = 91268.27
$ws.Range("I55").Value = 166983.17
$ws.Range("J55").Value = 410.4
$ws.Range("K55").Value = 166983.17
$ws.Range("L55").Value = 410.4
$ws.Range("M55").Value = -166769.17
$ws.Range("N55").Value = -838.4
$ws.Range("H60").Value = 3845.9092
$ws.Range("I60").Value = 4502.4287
$ws.Range("J60").Value = 2697
$ws.Range("K60").Value = 13507.2861
$ws.Range("L60").Value = 8091
$ws.Range("M60").Value = -13023.2861
$ws.Range("N60").Value = -9059
$ws.Range("H135").Value = 7898.3125
$ws.Range("I135").Value = 8902.643
$ws.Range("J135").Value = 868
$ws.Range("K135").Value = 80123.787
$ws.Range("L135").Value = 7812
$ws.Range("M135").Value = -77588.787
$ws.Range("N135").Value = -12882
$ws.Range("H141").Value = 2616.389
$ws.Range("I141").Value = 2120.8333
$ws.Range("J141").Value = 3607.5
$ws.Range("K141").Value = 6362.499899999999
$ws.Range("L141").Value = 10822.5
$ws.Range("M141").Value = -1182.499899999999
$ws.Range("N141").Value = -21182.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1085.6364
$ws.Range("I45").Value = 1044.2
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 1044.2
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -667.2
$ws.Range("N45").Value = -2254
$ws.Range("H74").Value = 5650.0356
$ws.Range("I74").Value = 1454.25
$ws.Range("J74").Value = 30824.75
$ws.Range("K74").Value = 1454.25
$ws.Range("L74").Value = 30824.75
$ws.Range("M74").Value = -580.25
$ws.Range("N74").Value = -32572.75
$ws.Range("H77").Value = 5650.0356
$ws.Range("I77").Value = 1454.25
$ws.Range("J77").Value = 30824.75
$ws.Range("K77").Value = 7271.25
$ws.Range("L77").Value = 154123.75
$ws.Range("M77").Value = -2903.25
$ws.Range("N77").Value = -162859.75
$ws.Range("H132").Value = 1526.9387
$ws.Range("I132").Value = 1241.8864
$ws.Range("J132").Value = 4035.4
$ws.Range("K132").Value = 3725.6592
$ws.Range("L132").Value = 12106.2
$ws.Range("M132").Value = -1195.6592
$ws.Range("N132").Value = -17166.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 49999.668
$ws.Range("J76").Value = 49999.668
$ws.Range("L76").Value = 49999.668
$ws.Range("N76").Value = -50629.668
$ws.Range("H79").Value = 49999.668
$ws.Range("J79").Value = 49999.668
$ws.Range("L79").Value = 49999.668
$ws.Range("N79").Value = -52183.668
$ws.Range("H134").Value = 3552
$ws.Range("I134").Value = 2539.2307
$ws.Range("J134").Value = 4564.769
$ws.Range("K134").Value = 7617.6921
$ws.Range("L134").Value = 13694.307
$ws.Range("M134").Value = -5082.6921
$ws.Range("N134").Value = -18764.307

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1578.4445
$ws.Range("I31").Value = 1473.2
$ws.Range("K31").Value = 1473.2
$ws.Range("M31").Value = -1178.2
$ws.Range("H34").Value = 1578.4445
$ws.Range("I34").Value = 1473.2
$ws.Range("K34").Value = 1473.2
$ws.Range("M34").Value = -1271.2
$ws.Range("H58").Value = 2790.2727
$ws.Range("I58").Value = 1956.9
$ws.Range("J58").Value = 3484.75
$ws.Range("K58").Value = 1956.9
$ws.Range("L58").Value = 3484.75
$ws.Range("M58").Value = -1753.9
$ws.Range("N58").Value = -3890.75
$ws.Range("H75").Value = 19900
$ws.Range("J75").Value = 19900
$ws.Range("L75").Value = 19900
$ws.Range("N75").Value = -21896
$ws.Range("H78").Value = 19900
$ws.Range("J78").Value = 19900
$ws.Range("L78").Value = 59700
$ws.Range("N78").Value = -69684
$ws.Range("H136").Value = 2790.2727
$ws.Range("I136").Value = 1956.9
$ws.Range("J136").Value = 3484.75
$ws.Range("K136").Value = 5870.700000000001
$ws.Range("L136").Value = 10454.25
$ws.Range("M136").Value = -3320.700000000001
$ws.Range("N136").Value = -15554.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2211.4546
$ws.Range("I5").Value = 1935.5
$ws.Range("K5").Value = 5806.5
$ws.Range("M5").Value = -5694.5
$ws.Range("H128").Value = 239096.75
$ws.Range("I128").Value = 239096.75
$ws.Range("K128").Value = 717290.25
$ws.Range("M128").Value = -712310.25
$ws.Range("H131").Value = 1657.1702
$ws.Range("I131").Value = 502.8
$ws.Range("J131").Value = 1969.1621
$ws.Range("K131").Value = 1508.4
$ws.Range("L131").Value = 5907.4863
$ws.Range("M131").Value = 3531.6
$ws.Range("N131").Value = -15987.4863
$ws.Range("H135").Value = 2211.4546
$ws.Range("I135").Value = 1935.5
$ws.Range("K135").Value = 17419.5
$ws.Range("M135").Value = -14884.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 23040.666
$ws.Range("I57").Value = 25000
$ws.Range("J57").Value = 22061
$ws.Range("K57").Value = 25000
$ws.Range("L57").Value = 22061
$ws.Range("M57").Value = -24180
$ws.Range("N57").Value = -23701
$ws.Range("H80").Value = 12760
$ws.Range("I80").Value = 15075
$ws.Range("K80").Value = 15075
$ws.Range("M80").Value = -14077
$ws.Range("H83").Value = 12760
$ws.Range("I83").Value = 15075
$ws.Range("K83").Value = 75375
$ws.Range("M83").Value = -70383
$ws.Range("H97").Value = 53233.105
$ws.Range("I97").Value = 72028.28999999999
$ws.Range("J97").Value = 606.6
$ws.Range("K97").Value = 72028.28999999999
$ws.Range("L97").Value = 606.6
$ws.Range("M97").Value = -71532.28999999999
$ws.Range("N97").Value = -1598.6
$ws.Range("H100").Value = 30000
$ws.Range("J100").Value = 30000
$ws.Range("L100").Value = 30000
$ws.Range("N100").Value = -32164

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 33449.75
$ws.Range("I76").Value = 7000
$ws.Range("J76").Value = 59899.5
$ws.Range("K76").Value = 7000
$ws.Range("L76").Value = 59899.5
$ws.Range("M76").Value = -6662
$ws.Range("N76").Value = -60575.5
$ws.Range("H79").Value = 33449.75
$ws.Range("I79").Value = 7000
$ws.Range("J79").Value = 59899.5
$ws.Range("K79").Value = 7000
$ws.Range("L79").Value = 59899.5
$ws.Range("M79").Value = -5830
$ws.Range("N79").Value = -62239.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 14150
$ws.Range("J45").Value = 14150
$ws.Range("L45").Value = 14150
$ws.Range("N45").Value = -15132
$ws.Range("H62").Value = 18897.5
$ws.Range("I62").Value = 31573
$ws.Range("K62").Value = 31573
$ws.Range("M62").Value = -30949
$ws.Range("H65").Value = 18897.5
$ws.Range("I65").Value = 31573
$ws.Range("K65").Value = 157865
$ws.Range("M65").Value = -154745
$ws.Range("H70").Value = 50000
$ws.Range("J70").Value = 50000
$ws.Range("L70").Value = 50000
$ws.Range("N70").Value = -50630
$ws.Range("H73").Value = 50000
$ws.Range("J73").Value = 50000
$ws.Range("L73").Value = 50000
$ws.Range("N73").Value = -52184
$ws.Range("H96").Value = 2524.75
$ws.Range("I96").Value = 2479.6
$ws.Range("J96").Value = 2600
$ws.Range("K96").Value = 2479.6
$ws.Range("L96").Value = 2600
$ws.Range("M96").Value = -1106.6
$ws.Range("N96").Value = -5346
$ws.Range("H107").Value = 487.36365
$ws.Range("I107").Value = 472.25
$ws.Range("J107").Value = 527.6667
$ws.Range("K107").Value = 1416.75
$ws.Range("L107").Value = 1583.0001
$ws.Range("M107").Value = 503.25
$ws.Range("N107").Value = -5423.0001
$ws.Range("H112").Value = 30667
$ws.Range("J112").Value = 30667
$ws.Range("L112").Value = 30667
$ws.Range("N112").Value = -33621
$ws.Range("H136").Value = 13932257
$ws.Range("I136").Value = 22289688
$ws.Range("J136").Value = 3203.2222
$ws.Range("K136").Value = 66869064
$ws.Range("L136").Value = 9609.6666
$ws.Range("M136").Value = -66866514
$ws.Range("N136").Value = -14709.6666
